# Insert a new weekly record at row 10 (shifting existing rows 10-37 down
# to rows 11-38), matching the "Fruta / hortaliza, semanal" update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift rows 10..37 down to 11..38, inserting a fresh blank row 10.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly observation.
$ws.Cells.Item(10, 1).Value  = 11
$ws.Cells.Item(10, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(10, 3).Value  = "Bíobío"
$ws.Cells.Item(10, 4).Value  = 44838
$ws.Cells.Item(10, 5).Value  = 8
$ws.Cells.Item(10, 6).Value  = 100112026
$ws.Cells.Item(10, 7).Value  = "Haba"
$ws.Cells.Item(10, 8).Value  = "Sin especificar"
$ws.Cells.Item(10, 9).Value  = "Primera"
$ws.Cells.Item(10, 10).Value = 100
$ws.Cells.Item(10, 11).Value = 11000
$ws.Cells.Item(10, 12).Value = 12000
$ws.Cells.Item(10, 13).Value = 11500
$ws.Cells.Item(10, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(10, 15).Value = "Región Metropolitana"
$ws.Cells.Item(10, 16).Value = 460
$ws.Cells.Item(10, 17).Value = 25
$ws.Cells.Item(10, 18).Value = "Hortaliza"
